$d = $word.ActiveDocument

# Anchor 1: the end of the "LOM3058: Química de Polímeros (Requisito fraco)"
# paragraph -- the blank paragraph, "Ver no Jupiter..." paragraph and the
# copyright/footer paragraph that immediately follow it are being removed.
$anchor1 = $d.Content
$anchor1.Find.Execute("LOM3058: Química de Polímeros (Requisito fraco)") | Out-Null
$anchor1.Expand(4) | Out-Null   # wdParagraph: expand to the whole paragraph
$deleteStart = $anchor1.End

# Anchor 2: the end of the copyright/footer paragraph ("... Powered by
# Jekyll and Github pages. Original theme under Creative Commons
# Attribution"). Everything from $deleteStart up to (but not including)
# this point -- i.e. the blank paragraph, the "Ver no Jupiter..."
# paragraph, and this footer paragraph itself, including its paragraph
# mark -- gets removed. The blank paragraph that originally followed the
# footer, and the page-break paragraph after it, are left untouched.
$anchor2 = $d.Content
$anchor2.Find.Execute("Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution") | Out-Null
$anchor2.Expand(4) | Out-Null
$deleteEnd = $anchor2.End

$d.Range($deleteStart, $deleteEnd).Delete()
